$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header cells (I1 "I0" and J1 "IF"), matching the style
# already used by the other header cells (e.g. H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data rows 2..34 for the new columns.
# Column I (I0) is always 1; column J (IF) mirrors column H for each row.
$lastRow = 34
for ($r = 2; $r -le $lastRow; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
